$wb = $excel.ActiveWorkbook

$sheetNames = @("Withdraw History", "Deposit History", "Transfer History", "Absolute History")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "Comments"
}

$wb.Worksheets.Item("Deposit History").Activate()
$wb.Worksheets.Item("Deposit History").Range("E1").Select()

$wb.Worksheets.Item("Transfer History").Activate()
$wb.Worksheets.Item("Transfer History").Range("E1").Select()

$wb.Worksheets.Item("Absolute History").Activate()
$wb.Worksheets.Item("Absolute History").Range("E5").Select()

$wb.Worksheets.Item("Withdraw History").Activate()
$wb.Worksheets.Item("Withdraw History").Range("E1").Select()
